$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294, shifting existing rows 294..399 down to 295..400
$ws.Rows(294).Insert()

# Populate the newly inserted row 294 with its data
$ws.Range("A294").Value = 10
$ws.Range("B294").Value = "Vega Modelo de Temuco"
$ws.Range("C294").Value = "La Araucanía"
$ws.Range("D294").Value = 45093
$ws.Range("E294").Value = 9
$ws.Range("F294").Value = 100112052
$ws.Range("G294").Value = "Albahaca"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 45
$ws.Range("K294").Value = 6000
$ws.Range("L294").Value = 6000
$ws.Range("M294").Value = 6000
$ws.Range("N294").Value = "`$/paquete"
$ws.Range("O294").Value = "Región de Arica y Parinacota"
$ws.Range("P294").Value = 6000
$ws.Range("Q294").Value = 1
$ws.Range("R294").Value = "Hortaliza"
